$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("Recorded By") holds a comma-separated list of the users who
# touched each attendance record. Reverse the order of that list for every
# data row (row 1 is the header) so the most-recently-recorded-by name
# moves to the front.
$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $reversedParts = $parts[($parts.Count - 1)..0]
            $cell.Value2 = [string]::Join(", ", $reversedParts)
        }
    }
}
